$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '72.965.27'
$ws.Range("E2").Value = '  +1.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '4.029.24'
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.65'
$ws.Range("E5").Value = '  +11.94%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.55'
$ws.Range("E6").Value = '  +1.45%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.687'
$ws.Range("E7").Value = '  -2.87%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.763'
$ws.Range("E9").Value = '  +1.88%  '

$ws.Range("E10").Value = '  -0.95%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.78'
$ws.Range("E11").Value = '  +13.23%  '

$ws.Range("E12").Value = '  -1.38%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.02'
$ws.Range("E13").Value = '  +3.70%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.663.02'
$ws.Range("E14").Value = '  -0.05%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.029.69'
$ws.Range("E15").Value = '  +0.27%  '

$ws.Range("E16").Value = '  +6.24%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.33'
$ws.Range("E17").Value = '  +2.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.71'
$ws.Range("E18").Value = '  +1.13%  '

$ws.Range("E19").Value = '  -0.36%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.656.40'
$ws.Range("E20").Value = '  +0.83%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '440.36'
$ws.Range("E21").Value = '  +2.79%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.72'
$ws.Range("E22").Value = '  +12.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '97.43'
$ws.Range("E23").Value = '  -0.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.55'
$ws.Range("E24").Value = '  +2.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.31'
$ws.Range("E25").Value = '  +0.84%  '

$ws.Range("E26").Value = '  +20.68%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.54'
$ws.Range("E27").Value = '  +1.45%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.78'
$ws.Range("E28").Value = '  +0.79%  '

$ws.Range("E29").Value = '  +1.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.84'
$ws.Range("E30").Value = '  +0.23%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.00'
$ws.Range("E31").Value = '  +11.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.134'
$ws.Range("E32").Value = '  +2.95%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '13.65'
$ws.Range("E33").Value = '  +2.40%  '

$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '49.86'
$ws.Range("E34").Value = '  +11.09%  '

$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '686.62'
$ws.Range("E35").Value = '  +0.81%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '71.32'
$ws.Range("E36").Value = '  +8.68%  '

$ws.Range("E37").Value = '  -0.41%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0877'
$ws.Range("E38").Value = '  +5.99%  '

$ws.Range("E39").Value = '  +5.38%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.149'
$ws.Range("E40").Value = '  -0.78%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.22'
$ws.Range("E41").Value = '  +13.97%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.37'
$ws.Range("E42").Value = '  -2.10%  '

$ws.Range("E43").Value = '  -0.20%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0493'
$ws.Range("E44").Value = '  +1.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.151'
$ws.Range("E46").Value = '  +0.20%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.73'
$ws.Range("E47").Value = '  +0.53%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.37'
$ws.Range("E48").Value = '  -0.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.52'
$ws.Range("E49").Value = '  +8.71%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.06'
$ws.Range("E50").Value = '  +1.36%  '

$ws.Range("E51").Value = '  +9.08%  '
